$wb = $excel.ActiveWorkbook

# --- ip_address_list: D11 -> remove second "OP:" line, keep only the XG-X2900 line ---
$ws1 = $wb.Worksheets.Item("ip_address_list")
$ws1.Range("D11").Value = "XG-X2900:`t`t10.101.28.175"

# --- ip_adress_fav_list: D4 -> same change ---
$ws2 = $wb.Worksheets.Item("ip_adress_fav_list")
$ws2.Range("D4").Value = "XG-X2900:`t`t10.101.28.175"

# --- Settings: B3 0 -> 1, B6 0 -> 1, add new row 9 ---
$ws4 = $wb.Worksheets.Item("Settings")
$ws4.Range("B3").Value = 1
$ws4.Range("B6").Value = 1
$ws4.Range("A9").Value = "automaticky přesouvat upravené projekty na začátek"
$ws4.Range("B9").Value = 1
